$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.109201
$ws.Range("H2").Value = 0.327603
$ws.Range("I2").Value = 0.07562717345335074
$ws.Range("J2").Value = 0.07562717345335074
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 0.1497010664763333
$ws.Range("R2").Value = 1.347309598287
$ws.Range("S2").Value = 0.0008342156022542801
$ws.Range("T2").Value = 0.00083421560225428
$ws.Range("G3").Value = 0.109201
$ws.Range("H3").Value = 0.327603
$ws.Range("I3").Value = 0.07562717345335074
$ws.Range("J3").Value = 0.07562717345335074
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 10.14627161092433
$ws.Range("R3").Value = 91.316444498319
$ws.Range("S3").Value = 0.05654053295526036
$ws.Range("T3").Value = 0.05654053295526035
$ws.Range("G4").Value = 0.109201
$ws.Range("H4").Value = 0.327603
$ws.Range("I4").Value = 0.07562717345335074
$ws.Range("J4").Value = 0.07562717345335074
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 3.245288171284
$ws.Range("R4").Value = 29.207593541556
$ws.Range("S4").Value = 0.01808450727853948
$ws.Range("T4").Value = 0.01808450727853948
$ws.Range("G5").Value = 0.109201
$ws.Range("H5").Value = 0.327603
$ws.Range("I5").Value = 0.07562717345335074
$ws.Range("J5").Value = 0.07562717345335074
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.030133033141
$ws.Range("R5").Value = 0.271197298269
$ws.Range("S5").Value = 0.0001679176172966112
$ws.Range("T5").Value = 0.0001679176172966112
$ws.Range("G6").Value = 0.7328223333333334
$ws.Range("I6").Value = 0.5075162472274908
$ws.Range("J6").Value = 0.5075162472274908
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 1.004608793304778
$ws.Range("R6").Value = 9.041479139743
$ws.Range("S6").Value = 0.005598225512102028
$ws.Range("T6").Value = 0.005598225512102028
$ws.Range("G7").Value = 0.7328223333333334
$ws.Range("I7").Value = 0.5075162472274908
$ws.Range("J7").Value = 0.5075162472274908
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("Q7").Value = 68.08925226464345
$ws.Range("R7").Value = 612.803270381791
$ws.Range("S7").Value = 0.3794302734240907
$ws.Range("T7").Value = 0.3794302734240907
$ws.Range("G8").Value = 0.7328223333333334
$ws.Range("I8").Value = 0.5075162472274908
$ws.Range("J8").Value = 0.5075162472274908
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 21.77836878800933
$ws.Range("R8").Value = 196.005319092084
$ws.Range("S8").Value = 0.1213608924922204
$ws.Range("T8").Value = 0.1213608924922203
$ws.Range("G9").Value = 0.7328223333333334
$ws.Range("I9").Value = 0.5075162472274908
$ws.Range("J9").Value = 0.5075162472274908
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 0.2022157274823333
$ws.Range("R9").Value = 1.819941547341
$ws.Range("S9").Value = 0.001126855799077631
$ws.Range("T9").Value = 0.00112685579907763
$ws.Range("G10").Value = 0.498848
$ws.Range("H10").Value = 1.496544
$ws.Range("I10").Value = 0.3454772778899196
$ws.Range("J10").Value = 0.3454772778899197
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 0.6838589171306667
$ws.Range("R10").Value = 6.154730254176001
$ws.Range("S10").Value = 0.003810833094507771
$ws.Range("T10").Value = 0.003810833094507771
$ws.Range("G11").Value = 0.498848
$ws.Range("H11").Value = 1.496544
$ws.Range("I11").Value = 0.3454772778899196
$ws.Range("J11").Value = 0.3454772778899197
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 46.34982555623467
$ws.Range("R11").Value = 417.1484300061121
$ws.Range("S11").Value = 0.2582863873377141
$ws.Range("T11").Value = 0.2582863873377141
$ws.Range("G12").Value = 0.498848
$ws.Range("H12").Value = 1.496544
$ws.Range("I12").Value = 0.3454772778899196
$ws.Range("J12").Value = 0.3454772778899197
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 14.825006306432
$ws.Range("R12").Value = 133.425056757888
$ws.Range("S12").Value = 0.08261298236174454
$ws.Range("T12").Value = 0.08261298236174454
$ws.Range("G13").Value = 0.498848
$ws.Range("H13").Value = 1.496544
$ws.Range("I13").Value = 0.3454772778899196
$ws.Range("J13").Value = 0.3454772778899197
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 0.137652615968
$ws.Range("R13").Value = 1.238873543712
$ws.Range("S13").Value = 0.0007670750959531498
$ws.Range("T13").Value = 0.0007670750959531498
$ws.Range("G14").Value = 0.1030673333333333
$ws.Range("H14").Value = 0.309202
$ws.Range("I14").Value = 0.07137930142923891
$ws.Range("J14").Value = 0.07137930142923891
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 0.1412925680064444
$ws.Range("R14").Value = 1.271633112058
$ws.Range("S14").Value = 0.0007873588845286152
$ws.Range("T14").Value = 0.0007873588845286151
$ws.Range("G15").Value = 0.1030673333333333
$ws.Range("H15").Value = 0.309202
$ws.Range("I15").Value = 0.07137930142923891
$ws.Range("J15").Value = 0.07137930142923891
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 9.57636979710511
$ws.Range("R15").Value = 86.187328173946
$ws.Range("S15").Value = 0.05336473069792528
$ws.Range("T15").Value = 0.05336473069792527
$ws.Range("G16").Value = 0.1030673333333333
$ws.Range("H16").Value = 0.309202
$ws.Range("I16").Value = 0.07137930142923891
$ws.Range("J16").Value = 0.07137930142923891
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 3.063004896589333
$ws.Range("R16").Value = 27.567044069304
$ws.Range("S16").Value = 0.01706872592601095
$ws.Range("T16").Value = 0.01706872592601095
$ws.Range("G17").Value = 0.1030673333333333
$ws.Range("H17").Value = 0.309202
$ws.Range("I17").Value = 0.07137930142923891
$ws.Range("J17").Value = 0.07137930142923891
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 0.02844050302733333
$ws.Range("R17").Value = 0.255964527246
$ws.Range("S17").Value = 0.0001584859207740673
$ws.Range("T17").Value = 0.0001584859207740673
